$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "70-53=17"
$t.Cell(1,2).Range.Text = "80-13=67"
$t.Cell(1,3).Range.Text = "54-26=28"
$t.Cell(1,4).Range.Text = "43-15=28"
$t.Cell(1,5).Range.Text = "45-29=16"

$t.Cell(2,1).Range.Text = "65+18=83"
$t.Cell(2,2).Range.Text = "17-8=9"
$t.Cell(2,3).Range.Text = "38+3=41"
$t.Cell(2,4).Range.Text = "68+3=71"
$t.Cell(2,5).Range.Text = "8+47=55"

$t.Cell(3,1).Range.Text = "63-15=48"
$t.Cell(3,2).Range.Text = "82-36=46"
$t.Cell(3,3).Range.Text = "38+7=45"
$t.Cell(3,4).Range.Text = "25+6=31"
$t.Cell(3,5).Range.Text = "15+38=53"

$t.Cell(4,1).Range.Text = "94-58=36"
$t.Cell(4,2).Range.Text = "55-36=19"
$t.Cell(4,3).Range.Text = "47+15=62"
$t.Cell(4,4).Range.Text = "24-5=19"
$t.Cell(4,5).Range.Text = "92-8=84"

$t.Cell(5,1).Range.Text = "35+47=82"
$t.Cell(5,2).Range.Text = "9+39=48"
$t.Cell(5,3).Range.Text = "6+49=55"
$t.Cell(5,4).Range.Text = "95-68=27"
$t.Cell(5,5).Range.Text = "97-49=48"

$t.Cell(6,1).Range.Text = "19+22=41"
$t.Cell(6,2).Range.Text = "36+45=81"
$t.Cell(6,3).Range.Text = "8+63=71"
$t.Cell(6,4).Range.Text = "24+57=81"
$t.Cell(6,5).Range.Text = "31-19=12"

$t.Cell(7,1).Range.Text = "19+34=53"
$t.Cell(7,2).Range.Text = "69+5=74"
$t.Cell(7,3).Range.Text = "59+8=67"
$t.Cell(7,4).Range.Text = "83-34=49"
$t.Cell(7,5).Range.Text = "42-34=8"

$t.Cell(8,1).Range.Text = "13-9=4"
$t.Cell(8,2).Range.Text = "15+48=63"
$t.Cell(8,3).Range.Text = "76+17=93"
$t.Cell(8,4).Range.Text = "81-15=66"
$t.Cell(8,5).Range.Text = "9+48=57"

$t.Cell(9,1).Range.Text = "18+65=83"
$t.Cell(9,2).Range.Text = "65+26=91"
$t.Cell(9,3).Range.Text = "69+16=85"
$t.Cell(9,4).Range.Text = "16+76=92"
$t.Cell(9,5).Range.Text = "36+49=85"

$t.Cell(10,1).Range.Text = "62-15=47"
$t.Cell(10,2).Range.Text = "6+69=75"
$t.Cell(10,3).Range.Text = "55+9=64"
$t.Cell(10,4).Range.Text = "83-76=7"
$t.Cell(10,5).Range.Text = "72-54=18"

$t.Cell(11,1).Range.Text = "37+25=62"
$t.Cell(11,2).Range.Text = "9+34=43"
$t.Cell(11,3).Range.Text = "48+38=86"
$t.Cell(11,4).Range.Text = "5+39=44"
$t.Cell(11,5).Range.Text = "8+74=82"

$t.Cell(12,1).Range.Text = "26+8=34"
$t.Cell(12,2).Range.Text = "70-15=55"
$t.Cell(12,3).Range.Text = "66+16=82"
$t.Cell(12,4).Range.Text = "23-5=18"
$t.Cell(12,5).Range.Text = "71-28=43"

$t.Cell(13,1).Range.Text = "45+26=71"
$t.Cell(13,2).Range.Text = "70-6=64"
$t.Cell(13,3).Range.Text = "79+6=85"
$t.Cell(13,4).Range.Text = "22-15=7"
$t.Cell(13,5).Range.Text = "52-7=45"

$t.Cell(14,1).Range.Text = "56+19=75"
$t.Cell(14,2).Range.Text = "80-17=63"
$t.Cell(14,3).Range.Text = "92-86=6"
$t.Cell(14,4).Range.Text = "32-27=5"
$t.Cell(14,5).Range.Text = "9+83=92"

$t.Cell(15,1).Range.Text = "16+19=35"
$t.Cell(15,2).Range.Text = "85-16=69"
$t.Cell(15,3).Range.Text = "73-48=25"
$t.Cell(15,4).Range.Text = "40-27=13"
$t.Cell(15,5).Range.Text = "8+67=75"

$t.Cell(16,1).Range.Text = "32-17=15"
$t.Cell(16,2).Range.Text = "17+74=91"
$t.Cell(16,3).Range.Text = "7+9=16"
$t.Cell(16,4).Range.Text = "80-25=55"
$t.Cell(16,5).Range.Text = "90-62=28"

$t.Cell(17,1).Range.Text = "38+23=61"
$t.Cell(17,2).Range.Text = "59+2=61"
$t.Cell(17,3).Range.Text = "62-43=19"
$t.Cell(17,4).Range.Text = "46+39=85"
$t.Cell(17,5).Range.Text = "82-28=54"

$t.Cell(18,1).Range.Text = "53-39=14"
$t.Cell(18,2).Range.Text = "35+56=91"
$t.Cell(18,3).Range.Text = "51-44=7"
$t.Cell(18,4).Range.Text = "36+57=93"
$t.Cell(18,5).Range.Text = "19+27=46"

$t.Cell(19,1).Range.Text = "6+85=91"
$t.Cell(19,2).Range.Text = "37+55=92"
$t.Cell(19,3).Range.Text = "65-27=38"
$t.Cell(19,4).Range.Text = "93-27=66"
$t.Cell(19,5).Range.Text = "81-73=8"

$t.Cell(20,1).Range.Text = "40-3=37"
$t.Cell(20,2).Range.Text = "97-38=59"
$t.Cell(20,3).Range.Text = "69+5=74"
$t.Cell(20,4).Range.Text = "29+29=58"
$t.Cell(20,5).Range.Text = "34+18=52"
